$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Preço" column (E) to hold the
# advertiser's name. This shifts E:L -> F:M.
$ws.Range("E1").EntireColumn.Insert()

# Populate the newly inserted column.
$ws.Range("E1").Value = "Nome do anunciante"
$ws.Range("E2").Value = "Richard"
